$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H19").Value = 1741.8667
$ws.Range("I19").Value = 1912.8
$ws.Range("J19").Value = 1400
$ws.Range("K19").Value = 1912.8
$ws.Range("L19").Value = 1400
$ws.Range("M19").Value = -1737.8
$ws.Range("N19").Value = -1750
$ws.Range("H40").Value = 26998.125
$ws.Range("I40").Value = 33999.668
$ws.Range("J40").Value = 22797.2
$ws.Range("K40").Value = 33999.668
$ws.Range("L40").Value = 22797.2
$ws.Range("M40").Value = -33824.668
$ws.Range("N40").Value = -23147.2
$ws.Range("H62").Value = 4142.7334
$ws.Range("I62").Value = 3514.4
$ws.Range("K62").Value = 3514.4
$ws.Range("M62").Value = -2890.4
$ws.Range("H65").Value = 4142.7334
$ws.Range("I65").Value = 3514.4
$ws.Range("K65").Value = 17572
$ws.Range("M65").Value = -14452
$ws.Range("H70").Value = 8242.786
$ws.Range("I70").Value = 6850
$ws.Range("J70").Value = 8799.9
$ws.Range("K70").Value = 20550
$ws.Range("L70").Value = 26399.7
$ws.Range("M70").Value = -20280
$ws.Range("N70").Value = -26939.7
$ws.Range("H73").Value = 8242.786
$ws.Range("I73").Value = 6850
$ws.Range("J73").Value = 8799.9
$ws.Range("K73").Value = 20550
$ws.Range("L73").Value = 26399.7
$ws.Range("M73").Value = -19614
$ws.Range("N73").Value = -28271.7
$ws.Range("H116").Value = 71440200
$ws.Range("I116").Value = 125007600
$ws.Range("K116").Value = 125007600
$ws.Range("M116").Value = -125004158
$ws.Range("H131").Value = 3269.7778
$ws.Range("I131").Value = 775.5714
$ws.Range("J131").Value = 11999.5
$ws.Range("K131").Value = 2326.7142
$ws.Range("L131").Value = 35998.5
$ws.Range("M131").Value = 2713.2858
$ws.Range("N131").Value = -46078.5
$ws.Range("H138").Value = 5940.635
$ws.Range("I138").Value = 2507
$ws.Range("J138").Value = 6150.857
$ws.Range("K138").Value = 7521
$ws.Range("L138").Value = 18452.571
$ws.Range("M138").Value = -2381
$ws.Range("N138").Value = -28732.571
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H32").Value = 7689.25
$ws.Range("I32").Value = 3776.652
$ws.Range("K32").Value = 3776.652
$ws.Range("M32").Value = -3489.652
$ws.Range("H63").Value = 999.5
$ws.Range("J63").Value = 1099.5
$ws.Range("L63").Value = 1099.5
$ws.Range("N63").Value = -2471.5
$ws.Range("H66").Value = 999.5
$ws.Range("J66").Value = 1099.5
$ws.Range("L66").Value = 5497.5
$ws.Range("N66").Value = -12361.5
$ws.Range("H106").Value = 64000
$ws.Range("J106").Value = 64000
$ws.Range("L106").Value = 64000
$ws.Range("N106").Value = -66524
$ws.Range("H132").Value = 27890.074
$ws.Range("I132").Value = 40308
$ws.Range("J132").Value = 12367.667
$ws.Range("K132").Value = 120924
$ws.Range("L132").Value = 37103.001
$ws.Range("M132").Value = -118394
$ws.Range("N132").Value = -42163.001
$ws.Range("H135").Value = 0
$ws.Range("J135").Value = 0
$ws.Range("L135").ClearContents()
$ws.Range("N135").Value = 0
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H20").Value = 2909.4546
$ws.Range("I20").Value = 3259
$ws.Range("K20").Value = 3259
$ws.Range("M20").Value = -3012
$ws.Range("H105").Value = 2594.25
$ws.Range("I105").Value = 2540.8696
$ws.Range("J105").Value = 2839.8
$ws.Range("K105").Value = 2540.8696
$ws.Range("L105").Value = 2839.8
$ws.Range("M105").Value = -793.8696
$ws.Range("N105").Value = -6333.8
$ws.Range("H107").Value = 9199
$ws.Range("I107").Value = 9199
$ws.Range("K107").Value = 9199
$ws.Range("M107").Value = -7279
$ws.Range("H130").Value = 0
$ws.Range("J130").Value = 0
$ws.Range("L130").ClearContents()
$ws.Range("N130").Value = 0
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H31").Value = 29414188
$ws.Range("J31").Value = 5473.375
$ws.Range("L31").Value = 5473.375
$ws.Range("N31").Value = -6063.375
$ws.Range("H34").Value = 29414188
$ws.Range("J34").Value = 5473.375
$ws.Range("L34").Value = 5473.375
$ws.Range("N34").Value = -5877.375
$ws.Range("H105").Value = 1150.1177
$ws.Range("I105").Value = 936.75
$ws.Range("K105").Value = 936.75
$ws.Range("M105").Value = 810.25
$ws.Range("H115").Value = 49246.668
$ws.Range("J115").Value = 49246.668
$ws.Range("L115").Value = 49246.668
$ws.Range("N115").Value = -51596.668
$ws.Range("H132").Value = 33345236
$ws.Range("I132").Value = 49387410
$ws.Range("K132").Value = 148162230
$ws.Range("M132").Value = -148159700
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H2").Value = 122.2
$ws.Range("I2").Value = 11.6
$ws.Range("J2").Value = 232.8
$ws.Range("K2").Value = 69.59999999999999
$ws.Range("L2").Value = 1396.8
$ws.Range("M2").Value = 43.40000000000001
$ws.Range("N2").Value = -1622.8
$ws.Range("H33").Value = 319.8
$ws.Range("I33").Value = 415.57144
$ws.Range("K33").Value = 2493.42864
$ws.Range("M33").Value = -2210.42864
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H25").Value = 2014
$ws.Range("I25").Value = 2014
$ws.Range("K25").Value = 2014
$ws.Range("M25").Value = -1485
$ws.Range("H80").Value = 114458.445
$ws.Range("I80").Value = 128265.75
$ws.Range("J80").Value = 4000
$ws.Range("K80").Value = 128265.75
$ws.Range("L80").Value = 4000
$ws.Range("M80").Value = -127267.75
$ws.Range("N80").Value = -5996
$ws.Range("H83").Value = 114458.445
$ws.Range("I83").Value = 128265.75
$ws.Range("J83").Value = 4000
$ws.Range("K83").Value = 641328.75
$ws.Range("L83").Value = 20000
$ws.Range("M83").Value = -636336.75
$ws.Range("N83").Value = -29984
$ws.Range("H123").Value = 36853
$ws.Range("J123").Value = 36853
$ws.Range("L123").Value = 36853
$ws.Range("N123").Value = -41753
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H7").Value = 3628.963
$ws.Range("I7").Value = 2995.574
$ws.Range("J7").Value = 4895.7407
$ws.Range("K7").Value = 2995.574
$ws.Range("L7").Value = 4895.7407
$ws.Range("M7").Value = -2883.574
$ws.Range("N7").Value = -5119.7407
$ws.Range("H43").Value = 0
$ws.Range("J43").Value = 0
$ws.Range("L43").ClearContents()
$ws.Range("N43").Value = 0
$ws.Range("H68").Value = 3022.7568
$ws.Range("I68").Value = 2621.15
$ws.Range("J68").Value = 3495.2354
$ws.Range("K68").Value = 2621.15
$ws.Range("L68").Value = 3495.2354
$ws.Range("M68").Value = -1872.15
$ws.Range("N68").Value = -4993.2354
$ws.Range("H71").Value = 3022.7568
$ws.Range("I71").Value = 2621.15
$ws.Range("J71").Value = 3495.2354
$ws.Range("K71").Value = 13105.75
$ws.Range("L71").Value = 17476.177
$ws.Range("M71").Value = -9361.75
$ws.Range("N71").Value = -24964.177
$ws.Range("H112").Value = 29950
$ws.Range("J112").Value = 29950
$ws.Range("L112").Value = 29950
$ws.Range("N112").Value = -32904
$ws.Range("H122").Value = 3807.121
$ws.Range("I122").Value = 2220.9375
$ws.Range("K122").Value = 6662.8125
$ws.Range("M122").Value = -4212.8125
$ws.Range("H126").Value = 3628.963
$ws.Range("I126").Value = 2995.574
$ws.Range("J126").Value = 4895.7407
$ws.Range("K126").Value = 8986.722
$ws.Range("L126").Value = 14687.2221
$ws.Range("M126").Value = -6516.722
$ws.Range("N126").Value = -19627.2221
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H46").Value = 57385.285
$ws.Range("J46").Value = 57385.285
$ws.Range("L46").Value = 57385.285
$ws.Range("N46").Value = -57847.285
$ws.Range("H62").Value = 8158.8
$ws.Range("I62").Value = 10000
$ws.Range("J62").Value = 7698.5
$ws.Range("K62").Value = 10000
$ws.Range("L62").Value = 7698.5
$ws.Range("M62").Value = -9376
$ws.Range("N62").Value = -8946.5
$ws.Range("H65").Value = 8158.8
$ws.Range("I65").Value = 10000
$ws.Range("J65").Value = 7698.5
$ws.Range("K65").Value = 50000
$ws.Range("L65").Value = 38492.5
$ws.Range("M65").Value = -46880
$ws.Range("N65").Value = -44732.5
$ws.Range("H122").Value = 2974.7827
$ws.Range("I122").Value = 1741.8889
$ws.Range("J122").Value = 7413.2
$ws.Range("K122").Value = 5225.6667
$ws.Range("L122").Value = 22239.6
$ws.Range("M122").Value = -2775.6667
$ws.Range("N122").Value = -27139.6
$ws.Range("H134").Value = 57385.285
$ws.Range("J134").Value = 57385.285
$ws.Range("L134").Value = 172155.855
$ws.Range("N134").Value = -177225.855
$ws.Range("H136").Value = 2516.3447
$ws.Range("I136").Value = 2031.2273
$ws.Range("K136").Value = 6093.6819
$ws.Range("M136").Value = -3543.6819
$ws.Range("H138").Value = 140370.33
$ws.Range("J138").Value = 195555.5
$ws.Range("L138").Value = 195555.5
$ws.Range("N138").Value = -205835.5
